$d = $word.ActiveDocument

# 1. Remove the standalone "Meta description: ..." paragraph that follows
#    the title heading (its content is being relocated to the end of the
#    document, see step 2/3 below).
$d.Paragraphs(2).Range.Delete()

# 2. Insert a new bold heading-like paragraph ("Play Fruit Slot for Free -
#    Board-Style Layout and Personalized Symbol Selection") right before the
#    final paragraph (the old "Create a feature image..." image-prompt
#    paragraph).
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex).Range
$lastPara.InsertParagraphBefore()

$newPara = $d.Paragraphs($lastIndex).Range
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruit Slot for Free - Board-Style Layout and Personalized Symbol Selection</w:t></w:r></w:p>'
$newPara.InsertXML($newParaXml)

# 3. Replace the text of the final paragraph (formerly the image-generation
#    prompt, still italic) with the meta description copy.
$oldText = "Create a feature image that captures the fun and excitement of Fruit Slot! Your image should be in a cartoon style and should feature a happy Maya warrior wearing glasses, surrounded by colorful fruit symbols. Be creative and playful with your design, incorporating the game's Asian arcade theme. Consider including the ring pattern of symbols in your design, as well as some of the potential multipliers that players can win. Your image should be eye-catching and convey the game's unique twist on traditional slots."
$newText = "Experience unique board-style layout with personalized symbol selection and occasional multipliers. Play Fruit Slot for free now."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
